$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("snapshot")

# Update scraped_at timestamps in column K (rows 2-42) on the "snapshot" sheet
$ws1.Cells.Item(2, 11).Value = "2025-11-20T07:02:29.457744+00:00"
$ws1.Cells.Item(3, 11).Value = "2025-11-20T07:02:29.457784+00:00"
$ws1.Cells.Item(4, 11).Value = "2025-11-20T07:02:29.457807+00:00"
$ws1.Cells.Item(5, 11).Value = "2025-11-20T07:02:31.249288+00:00"
$ws1.Cells.Item(6, 11).Value = "2025-11-20T07:02:31.249316+00:00"
$ws1.Cells.Item(7, 11).Value = "2025-11-20T07:02:33.417184+00:00"
$ws1.Cells.Item(8, 11).Value = "2025-11-20T07:02:35.610758+00:00"
$ws1.Cells.Item(9, 11).Value = "2025-11-20T07:02:37.472431+00:00"
$ws1.Cells.Item(10, 11).Value = "2025-11-20T07:02:37.472461+00:00"
$ws1.Cells.Item(11, 11).Value = "2025-11-20T07:02:37.472481+00:00"
$ws1.Cells.Item(12, 11).Value = "2025-11-20T07:02:39.686287+00:00"
$ws1.Cells.Item(13, 11).Value = "2025-11-20T07:02:41.863254+00:00"
$ws1.Cells.Item(14, 11).Value = "2025-11-20T07:02:43.982007+00:00"
$ws1.Cells.Item(15, 11).Value = "2025-11-20T07:02:46.164684+00:00"
$ws1.Cells.Item(16, 11).Value = "2025-11-20T07:02:46.164710+00:00"
$ws1.Cells.Item(17, 11).Value = "2025-11-20T07:02:46.164726+00:00"
$ws1.Cells.Item(18, 11).Value = "2025-11-20T07:02:48.337397+00:00"
$ws1.Cells.Item(19, 11).Value = "2025-11-20T07:02:50.059858+00:00"
$ws1.Cells.Item(20, 11).Value = "2025-11-20T07:02:52.168678+00:00"
$ws1.Cells.Item(21, 11).Value = "2025-11-20T07:02:52.168709+00:00"
$ws1.Cells.Item(22, 11).Value = "2025-11-20T07:02:54.341930+00:00"
$ws1.Cells.Item(23, 11).Value = "2025-11-20T07:02:54.341960+00:00"
$ws1.Cells.Item(24, 11).Value = "2025-11-20T07:02:54.341978+00:00"
$ws1.Cells.Item(25, 11).Value = "2025-11-20T07:02:56.171514+00:00"
$ws1.Cells.Item(26, 11).Value = "2025-11-20T07:02:56.171543+00:00"
$ws1.Cells.Item(27, 11).Value = "2025-11-20T07:02:57.919886+00:00"
$ws1.Cells.Item(28, 11).Value = "2025-11-20T07:02:57.919913+00:00"
$ws1.Cells.Item(29, 11).Value = "2025-11-20T07:02:57.919931+00:00"
$ws1.Cells.Item(30, 11).Value = "2025-11-20T07:02:59.747695+00:00"
$ws1.Cells.Item(31, 11).Value = "2025-11-20T07:02:59.747726+00:00"
$ws1.Cells.Item(32, 11).Value = "2025-11-20T07:03:02.003831+00:00"
$ws1.Cells.Item(33, 11).Value = "2025-11-20T07:03:02.003859+00:00"
$ws1.Cells.Item(34, 11).Value = "2025-11-20T07:03:02.003884+00:00"
$ws1.Cells.Item(35, 11).Value = "2025-11-20T07:03:02.003901+00:00"
$ws1.Cells.Item(36, 11).Value = "2025-11-20T07:03:02.003917+00:00"
$ws1.Cells.Item(37, 11).Value = "2025-11-20T07:03:04.006091+00:00"
$ws1.Cells.Item(38, 11).Value = "2025-11-20T07:03:04.006122+00:00"
$ws1.Cells.Item(39, 11).Value = "2025-11-20T07:03:08.646488+00:00"
$ws1.Cells.Item(40, 11).Value = "2025-11-20T07:03:08.646518+00:00"
$ws1.Cells.Item(41, 11).Value = "2025-11-20T07:03:08.646536+00:00"
$ws1.Cells.Item(42, 11).Value = "2025-11-20T07:03:10.557541+00:00"

# Remove the two new_injured rows (2 and 3), leaving only the header row
$ws3 = $wb.Worksheets.Item("new_injured")
$ws3.Rows("2:3").Delete()
